$wb = $excel.ActiveWorkbook

# --- Rename "Inspectare" sheet to "Requirements" ---
$wsReq = $wb.Worksheets.Item("Inspectare")
$wsReq.Name = "Requirements"

$wsArch = $wb.Worksheets.Item("Arhitectura")
$wsCode = $wb.Worksheets.Item("Cod Sursa")

# --- Fill in the "Observation" (column D) text cells in the exact order the ---
# --- original author typed them, so new shared-string table entries line up ---
$wsArch.Range("D5").Value  = "Chelnerul si clientul au acelasi terminal, se poate da place order & order served"
$wsArch.Range("D6").Value  = "Not specified (order 0 elements?)"
$wsArch.Range("D7").Value  = "Kitchen GUI nu e conectat la Controller; OrdersGUI nu are controller ci doar service"
$wsArch.Range("D10").Value = "Majoritatea da, dar nu este definita relatia dintre KitchenGUI si KitchenGUIController"
$wsArch.Range("D3").Value  = "Se reflecta arhitectura stratificata, model, repo, service, constroller, gui"
$wsArch.Range("D2").Value  = "In mare proiectul este organizat"
$wsArch.Range("D4").Value  = "Toate cerintele sunt respectate si au corespondent in diagrama de clase"
$wsArch.Range("D9").Value  = "Nume metode ambigue & lipsa, controller implementeaza actiuni ale service-ului"
$wsArch.Range("D11").Value = "Clasele principale corespund entitatilor din specificatii"
$wsArch.Range("D8").Value  = "PaymentAlert provoaca confuzie & MenuDataModel nu duce cu gandul la entitate (alternativa: Order)"

# --- "Cod Sursa" sheet gets its own new observation ---
$wsCode.Range("C9").Value = "NA?"

# --- Yes/NA answers (column C) on the Arhitectura sheet ---
$wsArch.Range("C2").Value  = "Yes"
$wsArch.Range("C3").Value  = "Yes"
$wsArch.Range("C4").Value  = "Yes"
$wsArch.Range("C5").Value  = "NA"
$wsArch.Range("C6").Value  = "NA"
$wsArch.Range("C7").Value  = "NA"
$wsArch.Range("C8").Value  = "NA"
$wsArch.Range("C9").Value  = "NA"
$wsArch.Range("C10").Value = "NA"
$wsArch.Range("C11").Value = "Yes"

# --- Wrap the long observation text so it matches the taller rows ---
$wsArch.Range("D3").WrapText = $true
$wsArch.Range("D4").WrapText = $true
$wsArch.Range("D5").WrapText = $true
$wsArch.Range("D7").WrapText = $true
$wsArch.Range("D8").WrapText = $true
$wsArch.Range("D9").WrapText = $true
$wsArch.Range("D10").WrapText = $true
$wsArch.Range("D11").WrapText = $true

# --- Row heights grew to fit the new wrapped text ---
$wsArch.Rows.Item(3).RowHeight = 45.75
$wsArch.Rows.Item(4).RowHeight = 42.75
$wsArch.Rows.Item(7).RowHeight = 48
$wsArch.Rows.Item(8).RowHeight = 60
$wsArch.Rows.Item(9).RowHeight = 42
$wsArch.Rows.Item(10).RowHeight = 45.75

# --- Selections left behind by the reviewer while editing ---
$wsReq.Activate()
$wsReq.Range("B11").Select()

$wsArch.Activate()
$wsArch.Range("D9").Select()
